$d = $word.ActiveDocument

$d.Content.Find.Execute("51×47=2397", $true, $false, $false, $false, $false, $true, 1, $false, "23×42=966", 2)
$d.Content.Find.Execute("59×86=5074", $true, $false, $false, $false, $false, $true, 1, $false, "44×51=2244", 2)
$d.Content.Find.Execute("53×36=1908", $true, $false, $false, $false, $false, $true, 1, $false, "98×55=5390", 2)
$d.Content.Find.Execute("16×70=1120", $true, $false, $false, $false, $false, $true, 1, $false, "45×50=2250", 2)
$d.Content.Find.Execute("91×65=5915", $true, $false, $false, $false, $false, $true, 1, $false, "47×78=3666", 2)
$d.Content.Find.Execute("69×61=4209", $true, $false, $false, $false, $false, $true, 1, $false, "17×60=1020", 2)
$d.Content.Find.Execute("98×23=2254", $true, $false, $false, $false, $false, $true, 1, $false, "53×35=1855", 2)
$d.Content.Find.Execute("66×63=4158", $true, $false, $false, $false, $false, $true, 1, $false, "33×77=2541", 2)
$d.Content.Find.Execute("39×44=1716", $true, $false, $false, $false, $false, $true, 1, $false, "34×41=1394", 2)
$d.Content.Find.Execute("94×72=6768", $true, $false, $false, $false, $false, $true, 1, $false, "98×27=2646", 2)
$d.Content.Find.Execute("51×11=561", $true, $false, $false, $false, $false, $true, 1, $false, "98×26=2548", 2)
$d.Content.Find.Execute("27×64=1728", $true, $false, $false, $false, $false, $true, 1, $false, "64×29=1856", 2)
$d.Content.Find.Execute("59×99=5841", $true, $false, $false, $false, $false, $true, 1, $false, "87×52=4524", 2)
$d.Content.Find.Execute("23×36=828", $true, $false, $false, $false, $false, $true, 1, $false, "59×58=3422", 2)
$d.Content.Find.Execute("52×88=4576", $true, $false, $false, $false, $false, $true, 1, $false, "87×81=7047", 2)
$d.Content.Find.Execute("56×48=2688", $true, $false, $false, $false, $false, $true, 1, $false, "83×56=4648", 2)
$d.Content.Find.Execute("53×47=2491", $true, $false, $false, $false, $false, $true, 1, $false, "20×82=1640", 2)
$d.Content.Find.Execute("24×92=2208", $true, $false, $false, $false, $false, $true, 1, $false, "64×55=3520", 2)
$d.Content.Find.Execute("19×31=589", $true, $false, $false, $false, $false, $true, 1, $false, "99×30=2970", 2)
$d.Content.Find.Execute("86×83=7138", $true, $false, $false, $false, $false, $true, 1, $false, "60×43=2580", 2)
$d.Content.Find.Execute("14×30=420", $true, $false, $false, $false, $false, $true, 1, $false, "33×47=1551", 2)
$d.Content.Find.Execute("75×66=4950", $true, $false, $false, $false, $false, $true, 1, $false, "25×93=2325", 2)
$d.Content.Find.Execute("40×28=1120", $true, $false, $false, $false, $false, $true, 1, $false, "47×34=1598", 2)
$d.Content.Find.Execute("69×16=1104", $true, $false, $false, $false, $false, $true, 1, $false, "19×75=1425", 2)
$d.Content.Find.Execute("71×74=5254", $true, $false, $false, $false, $false, $true, 1, $false, "87×51=4437", 2)
